# Correct swapped sex values for participants P14 and P15 (rows 26-29)
# and update the active selection, as part of starting the DXA/US analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# P14 (rows 26-27): was recorded as "f", should be "m"
$ws.Range("B26").Value = "m"
$ws.Range("B27").Value = "m"

# P15 (rows 28-29): was recorded as "m", should be "f"
$ws.Range("B28").Value = "f"
$ws.Range("B29").Value = "f"

# Update the current selection/view on the sheet
$ws.Range("B30").Select()
